# Add to cart feature - 18 Sept
# Insert a new MED_ID column at the front of the inventory table,
# populate it with medicine IDs, and fix the "PADADOL STRIP 20" typo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:E to B:F to make room for the new MED_ID column.
$ws.Columns("A").Insert()

# New MED_ID header + values.
$ws.Range("A1").Value = "MED_ID"
$ws.Range("A2").Value = 100
$ws.Range("A3").Value = 101
$ws.Range("A4").Value = 102
$ws.Range("A5").Value = 103

# Fix typo'd medicine name (was "PADADOL STRIP 20").
$ws.Range("B3").Value = "PANADOL STRIP 20"

# Match the author's final cell selection.
$ws.Range("B4").Select()
